$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the text-formatted "Price"/"Volume" columns must stay plain text even when
# the new value looks numeric, so pre-format those specific cells as Text first.

$ws.Range("D2").Value = "41.613.23"
$ws.Range("D3").Value = "2.254.25"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.08"
$ws.Range("E5").Value = "  +2.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.12"
$ws.Range("E6").Value = "  +3.99%  "
$ws.Range("E7").Value = "  +3.33%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.06"
$ws.Range("E10").Value = "  +3.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.81"
$ws.Range("E11").Value = "  +2.74%  "
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("E14").Value = "  +2.75%  "
$ws.Range("D15").Value = "2.603.45"
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.15"
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("D17").Value = "2.259.48"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D19").Value = "41.529.88"
$ws.Range("E19").Value = "  +3.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.28"
$ws.Range("E20").Value = "  +8.98%  "
$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.88"
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.53"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.00"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +5.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.87"
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.47"
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.05"
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.18"
$ws.Range("E32").Value = "  +6.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.13"
$ws.Range("E34").Value = "  +3.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0739"
$ws.Range("E35").Value = "  +3.52%  "
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.55"
$ws.Range("E39").Value = "  +5.54%  "
$ws.Range("E40").Value = "  +2.74%  "
$ws.Range("E41").Value = "  +2.55%  "
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("D43").Value = "2.051.04"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.51"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.16"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.85"
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("E48").Value = "  +6.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.52"
$ws.Range("E49").Value = "  +3.25%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.51"
$ws.Range("E50").Value = "  +7.01%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.15"
$ws.Range("E51").Value = "  +2.30%  "
